# İş Takip Güncellemesi - 25.11.2025 15:36:53
# Shifts İş Takip Listesi (rows 2-10, 33-122) J/K date-text cells back by one day,
# re-labels a handful of L-column status cells on the Güncelleme-sourced rows,
# and shifts the matching date-text cells (I/J/N/P) on the Güncelleme sheet back by one day.
#
# NOTE: every date in both sheets is stored as literal TEXT ("yyyy-mm-dd"), not a real
# Excel date serial (t="str" cells, no shared date format). Plain Range.Value assignment
# of a "yyyy-mm-dd"-looking string gets auto-coerced by Excel into a date serial number,
# so each date write below forces the cell to Text format first (and restores the Normal
# style afterwards so no stray formatting is left behind) to preserve the original text type.

$wb = $excel.ActiveWorkbook

# --- Sheet: İş Takip Listesi ---
$ws1 = $wb.Worksheets.Item("İş Takip Listesi")

# J (İŞE BAŞLAMA/YER TESLİMİ) / K (İHALE BİTİŞ TARİHİ) - each decremented by 1 day
$dateChanges1 = @(
    @{ Addr = "J2"; Value = "2025-08-13" },
    @{ Addr = "K2"; Value = "2026-01-13" },
    @{ Addr = "J3"; Value = "2025-08-13" },
    @{ Addr = "K3"; Value = "2026-01-13" },
    @{ Addr = "J4"; Value = "2025-08-13" },
    @{ Addr = "K4"; Value = "2026-01-13" },
    @{ Addr = "J5"; Value = "2025-08-13" },
    @{ Addr = "K5"; Value = "2026-01-13" },
    @{ Addr = "J6"; Value = "2025-08-13" },
    @{ Addr = "K6"; Value = "2026-01-13" },
    @{ Addr = "J7"; Value = "2025-08-13" },
    @{ Addr = "K7"; Value = "2026-01-13" },
    @{ Addr = "J8"; Value = "2025-08-13" },
    @{ Addr = "K8"; Value = "2026-01-13" },
    @{ Addr = "J9"; Value = "2025-08-13" },
    @{ Addr = "K9"; Value = "2026-01-13" },
    @{ Addr = "J10"; Value = "2025-08-13" },
    @{ Addr = "K10"; Value = "2026-01-13" },
    @{ Addr = "J33"; Value = "2025-08-15" },
    @{ Addr = "K33"; Value = "2026-01-15" },
    @{ Addr = "J34"; Value = "2025-08-15" },
    @{ Addr = "K34"; Value = "2026-01-15" },
    @{ Addr = "J35"; Value = "2025-08-15" },
    @{ Addr = "K35"; Value = "2026-01-15" },
    @{ Addr = "J36"; Value = "2025-08-15" },
    @{ Addr = "K36"; Value = "2026-01-15" },
    @{ Addr = "J37"; Value = "2025-08-15" },
    @{ Addr = "K37"; Value = "2026-01-15" },
    @{ Addr = "J38"; Value = "2025-08-15" },
    @{ Addr = "K38"; Value = "2026-01-15" },
    @{ Addr = "J39"; Value = "2025-08-15" },
    @{ Addr = "K39"; Value = "2026-01-15" },
    @{ Addr = "J40"; Value = "2025-08-15" },
    @{ Addr = "K40"; Value = "2026-01-15" },
    @{ Addr = "J41"; Value = "2025-08-15" },
    @{ Addr = "K41"; Value = "2026-01-15" },
    @{ Addr = "J42"; Value = "2025-08-15" },
    @{ Addr = "K42"; Value = "2026-01-15" },
    @{ Addr = "J43"; Value = "2025-08-15" },
    @{ Addr = "K43"; Value = "2026-01-15" },
    @{ Addr = "J44"; Value = "2025-08-15" },
    @{ Addr = "K44"; Value = "2026-01-15" },
    @{ Addr = "J45"; Value = "2025-08-15" },
    @{ Addr = "K45"; Value = "2026-01-15" },
    @{ Addr = "J46"; Value = "2025-08-15" },
    @{ Addr = "K46"; Value = "2026-01-15" },
    @{ Addr = "J47"; Value = "2025-08-15" },
    @{ Addr = "K47"; Value = "2026-01-15" },
    @{ Addr = "J48"; Value = "2025-08-15" },
    @{ Addr = "K48"; Value = "2026-01-15" },
    @{ Addr = "J49"; Value = "2025-08-15" },
    @{ Addr = "K49"; Value = "2026-01-15" },
    @{ Addr = "J50"; Value = "2025-08-15" },
    @{ Addr = "K50"; Value = "2026-01-15" },
    @{ Addr = "J51"; Value = "2025-08-15" },
    @{ Addr = "K51"; Value = "2026-01-15" },
    @{ Addr = "J52"; Value = "2025-08-15" },
    @{ Addr = "K52"; Value = "2026-01-15" },
    @{ Addr = "J53"; Value = "2025-08-15" },
    @{ Addr = "K53"; Value = "2026-01-15" },
    @{ Addr = "J54"; Value = "2025-08-15" },
    @{ Addr = "K54"; Value = "2026-01-15" },
    @{ Addr = "J55"; Value = "2025-08-15" },
    @{ Addr = "K55"; Value = "2026-01-15" },
    @{ Addr = "J56"; Value = "2025-08-15" },
    @{ Addr = "K56"; Value = "2026-01-15" },
    @{ Addr = "J57"; Value = "2025-08-15" },
    @{ Addr = "K57"; Value = "2026-01-15" },
    @{ Addr = "J58"; Value = "2025-08-15" },
    @{ Addr = "K58"; Value = "2026-01-15" },
    @{ Addr = "J59"; Value = "2025-08-15" },
    @{ Addr = "K59"; Value = "2026-01-15" },
    @{ Addr = "J60"; Value = "2025-08-15" },
    @{ Addr = "K60"; Value = "2026-01-15" },
    @{ Addr = "J61"; Value = "2025-08-15" },
    @{ Addr = "K61"; Value = "2026-01-15" },
    @{ Addr = "J62"; Value = "2025-08-15" },
    @{ Addr = "K62"; Value = "2026-01-15" },
    @{ Addr = "J63"; Value = "2025-08-15" },
    @{ Addr = "K63"; Value = "2026-01-15" },
    @{ Addr = "J64"; Value = "2025-08-15" },
    @{ Addr = "K64"; Value = "2026-01-15" },
    @{ Addr = "J65"; Value = "2025-08-15" },
    @{ Addr = "K65"; Value = "2026-01-15" },
    @{ Addr = "J66"; Value = "2025-08-15" },
    @{ Addr = "K66"; Value = "2026-01-15" },
    @{ Addr = "J67"; Value = "2025-08-15" },
    @{ Addr = "K67"; Value = "2026-01-15" },
    @{ Addr = "J68"; Value = "2025-08-15" },
    @{ Addr = "K68"; Value = "2026-01-15" },
    @{ Addr = "J69"; Value = "2025-08-15" },
    @{ Addr = "K69"; Value = "2026-01-15" },
    @{ Addr = "J70"; Value = "2025-08-15" },
    @{ Addr = "K70"; Value = "2026-01-15" },
    @{ Addr = "J71"; Value = "2025-08-15" },
    @{ Addr = "K71"; Value = "2026-01-15" },
    @{ Addr = "J72"; Value = "2025-08-15" },
    @{ Addr = "K72"; Value = "2026-01-15" },
    @{ Addr = "J73"; Value = "2025-08-15" },
    @{ Addr = "K73"; Value = "2026-01-15" },
    @{ Addr = "J74"; Value = "2025-08-15" },
    @{ Addr = "K74"; Value = "2026-01-15" },
    @{ Addr = "J75"; Value = "2025-08-15" },
    @{ Addr = "K75"; Value = "2026-01-15" },
    @{ Addr = "J76"; Value = "2025-08-15" },
    @{ Addr = "K76"; Value = "2026-01-15" },
    @{ Addr = "J77"; Value = "2025-08-15" },
    @{ Addr = "K77"; Value = "2026-01-15" },
    @{ Addr = "J78"; Value = "2025-08-15" },
    @{ Addr = "K78"; Value = "2026-01-15" },
    @{ Addr = "J79"; Value = "2025-08-15" },
    @{ Addr = "K79"; Value = "2026-01-15" },
    @{ Addr = "J80"; Value = "2025-08-15" },
    @{ Addr = "K80"; Value = "2026-01-15" },
    @{ Addr = "J81"; Value = "2025-08-15" },
    @{ Addr = "K81"; Value = "2026-01-15" },
    @{ Addr = "J82"; Value = "2025-08-15" },
    @{ Addr = "K82"; Value = "2026-01-15" },
    @{ Addr = "J83"; Value = "2025-08-15" },
    @{ Addr = "K83"; Value = "2026-01-15" },
    @{ Addr = "J84"; Value = "2025-08-15" },
    @{ Addr = "K84"; Value = "2026-01-15" },
    @{ Addr = "J85"; Value = "2025-08-15" },
    @{ Addr = "K85"; Value = "2026-01-15" },
    @{ Addr = "J86"; Value = "2025-08-15" },
    @{ Addr = "K86"; Value = "2026-01-15" },
    @{ Addr = "J87"; Value = "2025-08-15" },
    @{ Addr = "K87"; Value = "2026-01-15" },
    @{ Addr = "J88"; Value = "2025-08-15" },
    @{ Addr = "K88"; Value = "2026-01-15" },
    @{ Addr = "J89"; Value = "2025-08-15" },
    @{ Addr = "K89"; Value = "2026-01-15" },
    @{ Addr = "J90"; Value = "2025-08-15" },
    @{ Addr = "K90"; Value = "2026-01-15" },
    @{ Addr = "J91"; Value = "2025-08-15" },
    @{ Addr = "K91"; Value = "2026-01-15" },
    @{ Addr = "J92"; Value = "2025-08-15" },
    @{ Addr = "K92"; Value = "2026-01-15" },
    @{ Addr = "J93"; Value = "2025-08-15" },
    @{ Addr = "K93"; Value = "2026-01-15" },
    @{ Addr = "J94"; Value = "2025-08-15" },
    @{ Addr = "K94"; Value = "2026-01-15" },
    @{ Addr = "J95"; Value = "2024-06-13" },
    @{ Addr = "K95"; Value = "2025-08-07" },
    @{ Addr = "J96"; Value = "2024-06-13" },
    @{ Addr = "K96"; Value = "2025-08-07" },
    @{ Addr = "J97"; Value = "2024-06-13" },
    @{ Addr = "K97"; Value = "2025-08-07" },
    @{ Addr = "J98"; Value = "2024-06-13" },
    @{ Addr = "K98"; Value = "2025-08-07" },
    @{ Addr = "J99"; Value = "2024-06-13" },
    @{ Addr = "K99"; Value = "2025-08-07" },
    @{ Addr = "J100"; Value = "2024-06-13" },
    @{ Addr = "K100"; Value = "2025-08-07" },
    @{ Addr = "J101"; Value = "2024-06-13" },
    @{ Addr = "K101"; Value = "2025-08-07" },
    @{ Addr = "J102"; Value = "2024-06-13" },
    @{ Addr = "K102"; Value = "2025-08-07" },
    @{ Addr = "J103"; Value = "2024-06-13" },
    @{ Addr = "K103"; Value = "2025-08-07" },
    @{ Addr = "J104"; Value = "2024-06-13" },
    @{ Addr = "K104"; Value = "2025-08-07" },
    @{ Addr = "J105"; Value = "2024-06-13" },
    @{ Addr = "K105"; Value = "2025-08-07" },
    @{ Addr = "J106"; Value = "2024-06-13" },
    @{ Addr = "K106"; Value = "2025-08-07" },
    @{ Addr = "J107"; Value = "2024-06-13" },
    @{ Addr = "K107"; Value = "2025-08-07" },
    @{ Addr = "J108"; Value = "2024-06-13" },
    @{ Addr = "K108"; Value = "2025-08-07" },
    @{ Addr = "J109"; Value = "2024-06-13" },
    @{ Addr = "K109"; Value = "2025-08-07" },
    @{ Addr = "J110"; Value = "2024-06-13" },
    @{ Addr = "K110"; Value = "2025-08-07" },
    @{ Addr = "J111"; Value = "2024-06-13" },
    @{ Addr = "K111"; Value = "2025-08-07" },
    @{ Addr = "J112"; Value = "2024-06-13" },
    @{ Addr = "K112"; Value = "2025-08-07" },
    @{ Addr = "J113"; Value = "2024-06-13" },
    @{ Addr = "K113"; Value = "2025-08-07" },
    @{ Addr = "J114"; Value = "2024-06-13" },
    @{ Addr = "K114"; Value = "2025-08-07" },
    @{ Addr = "J115"; Value = "2024-06-13" },
    @{ Addr = "K115"; Value = "2025-08-07" },
    @{ Addr = "J116"; Value = "2024-06-13" },
    @{ Addr = "K116"; Value = "2025-08-07" },
    @{ Addr = "J117"; Value = "2024-06-13" },
    @{ Addr = "K117"; Value = "2025-08-07" },
    @{ Addr = "J118"; Value = "2024-06-13" },
    @{ Addr = "K118"; Value = "2025-08-07" },
    @{ Addr = "J119"; Value = "2024-06-13" },
    @{ Addr = "K119"; Value = "2025-08-07" },
    @{ Addr = "J120"; Value = "2024-06-13" },
    @{ Addr = "K120"; Value = "2025-08-07" },
    @{ Addr = "J121"; Value = "2024-06-13" },
    @{ Addr = "K121"; Value = "2025-08-07" },
    @{ Addr = "J122"; Value = "2024-06-13" },
    @{ Addr = "K122"; Value = "2025-08-07" }
)
foreach ($chg in $dateChanges1) {
    $cell = $ws1.Range($chg.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
    $cell.Style = "Normal"
}

# L (DURUMU) - plain text status relabels, no date coercion risk
$textChanges1 = @(
    @{ Addr = "L60"; Value = "KESİN ASKIDA" },
    @{ Addr = "L62"; Value = "KESİN ASKIDA" },
    @{ Addr = "L76"; Value = "KESİN ASKIDA" },
    @{ Addr = "L79"; Value = "KESİN ASKIDA" },
    @{ Addr = "L80"; Value = "KESİN ASKIDA" },
    @{ Addr = "L86"; Value = "KESİN ASKIDA" },
    @{ Addr = "L89"; Value = "KESİN ASKIDA" }
)
foreach ($chg in $textChanges1) {
    $ws1.Range($chg.Addr).Value = $chg.Value
}

# --- Sheet: Güncelleme ---
$ws2 = $wb.Worksheets.Item("Güncelleme")

# I/J/N/P date-text columns - each decremented by 1 day (blank cells left untouched)
$dateChanges2 = @(
    @{ Addr = "J2"; Value = "2024-09-19" },
    @{ Addr = "N2"; Value = "2025-05-23" },
    @{ Addr = "P2"; Value = "2025-08-09" },
    @{ Addr = "J3"; Value = "2024-12-21" },
    @{ Addr = "N3"; Value = "2025-09-10" },
    @{ Addr = "J4"; Value = "2024-10-25" },
    @{ Addr = "N4"; Value = "2025-04-17" },
    @{ Addr = "P4"; Value = "2025-07-12" },
    @{ Addr = "I5"; Value = "2025-04-20" },
    @{ Addr = "J6"; Value = "2025-12-01" },
    @{ Addr = "N6"; Value = "2025-08-20" },
    @{ Addr = "I7"; Value = "2024-12-21" },
    @{ Addr = "J7"; Value = "2024-12-21" },
    @{ Addr = "J8"; Value = "2024-12-09" },
    @{ Addr = "N8"; Value = "2025-05-10" },
    @{ Addr = "P8"; Value = "2025-06-12" },
    @{ Addr = "I9"; Value = "2025-08-06" },
    @{ Addr = "J9"; Value = "2025-01-22" },
    @{ Addr = "J10"; Value = "2024-11-20" },
    @{ Addr = "N10"; Value = "2025-08-31" },
    @{ Addr = "I11"; Value = "2025-05-28" },
    @{ Addr = "J11"; Value = "2025-01-03" },
    @{ Addr = "N11"; Value = "2025-09-20" },
    @{ Addr = "J12"; Value = "2024-12-01" },
    @{ Addr = "N12"; Value = "2025-08-10" },
    @{ Addr = "J13"; Value = "2025-01-29" },
    @{ Addr = "J14"; Value = "2025-11-27" },
    @{ Addr = "J15"; Value = "2025-02-17" },
    @{ Addr = "N15"; Value = "2025-09-07" },
    @{ Addr = "J16"; Value = "2024-10-16" },
    @{ Addr = "N16"; Value = "2025-03-26" },
    @{ Addr = "P16"; Value = "2025-06-12" },
    @{ Addr = "J17"; Value = "2024-11-01" },
    @{ Addr = "J18"; Value = "2025-04-10" },
    @{ Addr = "I19"; Value = "2025-05-29" },
    @{ Addr = "J19"; Value = "2025-02-17" },
    @{ Addr = "N19"; Value = "2025-09-27" },
    @{ Addr = "J20"; Value = "2025-01-29" },
    @{ Addr = "J21"; Value = "2024-11-22" },
    @{ Addr = "J22"; Value = "2024-11-22" },
    @{ Addr = "J23"; Value = "2025-01-30" },
    @{ Addr = "I24"; Value = "2025-07-27" },
    @{ Addr = "J25"; Value = "2024-12-25" },
    @{ Addr = "J27"; Value = "2025-03-17" },
    @{ Addr = "J28"; Value = "2025-01-13" },
    @{ Addr = "I29"; Value = "2025-04-04" },
    @{ Addr = "J29"; Value = "2025-01-30" }
)
foreach ($chg in $dateChanges2) {
    $cell = $ws2.Range($chg.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
    $cell.Style = "Normal"
}

